$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# A1: date serial 45310 -> 45311 (2024-01-19 -> 2024-01-20)
$ws.Range("A1").Value = 45311

# D14:D21 price updates
$ws.Range("D14").Value = 20.3
$ws.Range("D15").Value = 23.2
$ws.Range("D16").Value = 28.01
$ws.Range("D17").Value = 35.7
$ws.Range("D18").Value = 71.3
$ws.Range("D19").Value = 107
$ws.Range("D20").Value = 168.3
$ws.Range("D21").Value = 233
